$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.577.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.62%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.604.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.68%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.77%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("E8").Value = "  +0.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.54"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("E10").Value = "  +1.11%  "

$ws.Range("E11").Value = "  +1.64%  "

$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.063.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.487.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.87%  "

$ws.Range("E15").Value = "  +0.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.659.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.52%  "

$ws.Range("E17").Value = "  +0.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "341.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("E19").Value = "  +1.62%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.74%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("E23").Value = "  +1.57%  "

$ws.Range("E24").Value = "  +1.42%  "

$ws.Range("E25").Value = "  -1.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.994"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.46%  "

$ws.Range("E28").Value = "  +3.05%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +6.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.56%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.56%  "

$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.848"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.52%  "

$ws.Range("E37").Value = "  -0.80%  "

$ws.Range("E38").Value = "  +0.27%  "

$ws.Range("E39").Value = "  +0.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "273.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.45%  "

$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("E43").Value = "  -0.67%  "

$ws.Range("E44").Value = "  +0.17%  "

$ws.Range("E45").Value = "  +1.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.53%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0223"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.92%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.939.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.50%  "

$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.60%  "

$ws.Range("E51").Value = "  +2.10%  "
